$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-modified date for every data row.
# Update every data row (2 through 490) from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C490").Value = 45177
